$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated crypto price (column D) and 1h volume change (column E)
# values scraped by the GitHub Actions workflow.
#
# Several of the new column D values (e.g. "1.000", "5.322") are strings that
# look like plain numbers and would otherwise be auto-converted by Excel,
# losing the trailing/formatting digits (and column D elsewhere legitimately
# holds dotted-thousands text like "26.877.44"). For every column D cell we
# therefore switch the cell to a text number format before writing the value,
# then restore the default "Normal" cell style so no stray formatting is left
# behind.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '26.877.44'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.40%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.815.85'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  +0.13%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '308.88'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  +0.06%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.4667'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +0.66%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3685'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("E9").Value = '  +1.22%  '
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("E11").Value = '  -0.08%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '1.752.81'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +1.31%  '
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("E14").Value = '  +0.00%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.07062'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +0.25%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '91.55'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("E17").Value = '  +0.25%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.000008694'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("E19").Value = '  +0.04%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '14.71'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.61%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '26.913.89'
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '5.322'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.64%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '10.59'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.52%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.040.82'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +3.99%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '1.895'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.65%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '150.12'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.19%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '2.169'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +1.54%  '
$ws.Range("E28").Value = '  +0.69%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '5.324'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +1.90%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '115.75'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +1.38%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.08913'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +0.10%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.7670'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +1.19%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.163'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.08%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '4.504'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +1.61%  '
$ws.Range("E35").Value = '  +0.40%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +0.06%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '1.085'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -3.09%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.01960'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.91%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.05286'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +1.55%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.929'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +0.88%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '7.256'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +1.07%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.5321'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +2.00%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '2.348'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -1.45%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.1659'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.65%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '8.427'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -0.89%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.4923'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -1.66%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '10.43'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +1.57%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("E49").Value = '  +1.34%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '103.79'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -0.43%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.06286'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.05%  '
